# Fill in the empty references in the "Payment Table" design slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table

# Row 6 = report_no, column 4 = Format/Remark reference cell (currently "-")
$cell6 = $tbl.Cell(6, 4)
$cell6.Shape.TextFrame.TextRange.Text = "Reference Report Table"
$cell6.Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 1

# Row 7 = account_no, column 4 = Format/Remark reference cell (currently "-")
$cell7 = $tbl.Cell(7, 4)
$cell7.Shape.TextFrame.TextRange.Text = "Reference Account Table"
$cell7.Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 1
